# Applies the "Add files via upload" edit to Saldo_guide.xlsx
#  1) Shift the reference date in column G (rows 2-310) from 2024-04-08 (45390)
#     to 2024-04-09 (45391) for every data row.
#  2) Update a handful of rows whose Saldo Previsto (D), Vl. Projetado (E)
#     and/or Vl. Total (H) values were recalculated for the new date.
#  3) Move the active cell selection from O10 to N9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Saldo_guide")

# --- 1) Shift every reference date in column G (rows 2..310) by one day ---
$lastRow = 310
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $cell.Value2 = $cell.Value2 + 1
}

# --- 2) Row-specific updates to D (Saldo Previsto), E (Vl. Projetado) and H (Vl. Total) ---
$rowUpdates = @{
    2   = @{ D = 93589.64;   E = 0 }
    27  = @{ D = 4135.57;    E = 0;        H = 4135.57 }
    61  = @{ D = 4295.38;    E = 0 }
    71  = @{ D = 11447.65;                 H = 11447.65 }
    78  = @{ D = 603.99;     E = 0 }
    115 = @{ D = 943.85;     E = 0;        H = 943.85 }
    118 = @{ D = 4090.58;                  H = 4090.58 }
    119 = @{ D = 105457.38;                H = 105457.38 }
    120 = @{ D = 766.44;     E = 0;        H = 766.44 }
    125 = @{               E = -29162.65;  H = 703.2 }
    151 = @{ D = 18242.84;                 H = 18242.84 }
    167 = @{ D = 164848.27;  E = -158889.97; H = 5958.3 }
    175 = @{               E = -80.61;     H = 731.67 }
    184 = @{               E = -41.97;     H = 10869.48 }
    192 = @{               E = -887;       H = 47.92 }
    200 = @{ D = 1385.12;   E = -132.06;   H = 1253.06 }
    220 = @{               E = -64.33;     H = 578.74 }
    239 = @{               E = -73.92;     H = 848.3 }
    266 = @{               E = -228.42;    H = 330.91 }
    267 = @{ D = 2494.31;   E = -88.93;    H = 2405.38 }
    268 = @{               E = -316.34;    H = -304.89 }
    274 = @{               E = -150.61;    H = 47.66 }
    280 = @{               E = -6.08;      H = 102.5 }
}

foreach ($rowNum in $rowUpdates.Keys) {
    $vals = $rowUpdates[$rowNum]
    if ($vals.ContainsKey("D")) { $ws.Cells.Item($rowNum, 4).Value2 = $vals["D"] }
    if ($vals.ContainsKey("E")) { $ws.Cells.Item($rowNum, 5).Value2 = $vals["E"] }
    if ($vals.ContainsKey("H")) { $ws.Cells.Item($rowNum, 8).Value2 = $vals["H"] }
}

# --- 3) Move active cell selection from O10 to N9 ---
$ws.Range("N9").Select()
